# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# Business changes to the "Hoja1" Estado de Cuenta sheet:
#   - Valor Mora (E11) updated: 276000 -> 184000
#   - Cant. Periodos (F13) updated: 3 -> 2 (one period removed)
#   - The middle detail row (period 2504, row 17) is removed entirely;
#     the old last detail row (period 2503, row 18) shifts up to row 17
#     and its period value is updated to 2505.
#   - Signature rows shift up from 23/24 to 22/23 as a consequence of the
#     deleted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "2504" detail row (row 17) — the other rows (16, 18) shift up
# and keep their own per-cell styles, matching how Excel handles a row
# delete in the middle of the worker/period table.
$ws.Rows("17").Delete()

# After the delete, row 16 still holds the first worker/period entry and
# row 17 (previously row 18) holds the last one. Update the period labels
# for the two remaining rows.
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"

# Update the "Valor Mora" total and "Cant. Periodos" count.
$ws.Range("E11").Value = 184000
$ws.Range("F13").Value = 2
